$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value2 = 46081

# Row 3
$ws.Range("C3").Value2 = 46081

# Row 4
$ws.Range("C4").Value2 = 46081

# Row 5
$ws.Range("A5").Value2 = 'A 61558-2023'
$ws.Range("B5").Value2 = 45265
$ws.Range("C5").Value2 = 46081
$ws.Range("F5").Value2 = 'Övriga statliga verk och myndigheter'
$ws.Range("G5").Value2 = 1.5
$ws.Range("J5").Value2 = 0
$ws.Range("L5").Value2 = 1
$ws.Range("P5").Value2 = 1
$ws.Range("R5").Value2 = 'Ask'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 61558-2023 artfynd.xlsx", "A 61558-2023")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 61558-2023 karta.png", "A 61558-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 61558-2023 FSC-klagomål.docx", "A 61558-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 61558-2023 FSC-klagomål mail.docx", "A 61558-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 61558-2023 tillsynsbegäran.docx", "A 61558-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 61558-2023 tillsynsbegäran mail.docx", "A 61558-2023")'

# Row 6
$ws.Range("A6").Value2 = 'A 34341-2024'
$ws.Range("B6").Value2 = 45525
$ws.Range("C6").Value2 = 46081
$ws.Range("F6").Value2 = 'Övriga Aktiebolag'
$ws.Range("G6").Value2 = 14.4
$ws.Range("R6").Value2 = 'Desmeknopp'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 34341-2024 artfynd.xlsx", "A 34341-2024")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 34341-2024 karta.png", "A 34341-2024")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 34341-2024 FSC-klagomål.docx", "A 34341-2024")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 34341-2024 FSC-klagomål mail.docx", "A 34341-2024")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 34341-2024 tillsynsbegäran.docx", "A 34341-2024")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 34341-2024 tillsynsbegäran mail.docx", "A 34341-2024")'

# Row 7
$ws.Range("A7").Value2 = 'A 31213-2023'
$ws.Range("B7").Value2 = 45113
$ws.Range("C7").Value2 = 46081
$ws.Range("G7").Value2 = 6.5
$ws.Range("H7").Value2 = 0
$ws.Range("J7").Value2 = 1
$ws.Range("O7").Value2 = 1
$ws.Range("R7").Value2 = 'Skogsveronika'
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 31213-2023 artfynd.xlsx", "A 31213-2023")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 31213-2023 karta.png", "A 31213-2023")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 31213-2023 FSC-klagomål.docx", "A 31213-2023")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 31213-2023 FSC-klagomål mail.docx", "A 31213-2023")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 31213-2023 tillsynsbegäran.docx", "A 31213-2023")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 31213-2023 tillsynsbegäran mail.docx", "A 31213-2023")'

# Row 8
$ws.Range("A8").Value2 = 'A 13766-2023'
$ws.Range("B8").Value2 = 45007
$ws.Range("C8").Value2 = 46081
$ws.Range("G8").Value2 = 0.9
$ws.Range("R8").Value2 = 'Större vattensalamander'
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 13766-2023 artfynd.xlsx", "A 13766-2023")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 13766-2023 karta.png", "A 13766-2023")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 13766-2023 FSC-klagomål.docx", "A 13766-2023")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 13766-2023 FSC-klagomål mail.docx", "A 13766-2023")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 13766-2023 tillsynsbegäran.docx", "A 13766-2023")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 13766-2023 tillsynsbegäran mail.docx", "A 13766-2023")'

# Row 9
$ws.Range("A9").Value2 = 'A 60891-2024'
$ws.Range("B9").Value2 = 45644
$ws.Range("C9").Value2 = 46081
$ws.Range("F9").Value2 = ""
$ws.Range("G9").Value2 = 16.1
$ws.Range("H9").Value2 = 1
$ws.Range("L9").Value2 = 0
$ws.Range("O9").Value2 = 0
$ws.Range("P9").Value2 = 0
$ws.Range("R9").Value2 = 'Lövgroda'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/artfynd/A 60891-2024 artfynd.xlsx", "A 60891-2024")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/kartor/A 60891-2024 karta.png", "A 60891-2024")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomål/A 60891-2024 FSC-klagomål.docx", "A 60891-2024")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/klagomålsmail/A 60891-2024 FSC-klagomål mail.docx", "A 60891-2024")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsyn/A 60891-2024 tillsynsbegäran.docx", "A 60891-2024")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/tillsynsmail/A 60891-2024 tillsynsbegäran mail.docx", "A 60891-2024")'

# Row 10
$ws.Range("C10").Value2 = 46081

# Row 11
$ws.Range("C11").Value2 = 46081

# Row 12
$ws.Range("C12").Value2 = 46081

# Row 13
$ws.Range("C13").Value2 = 46081

# Row 14
$ws.Range("A14").Value2 = 'A 40417-2022'
$ws.Range("B14").Value2 = 44823
$ws.Range("C14").Value2 = 46081
$ws.Range("G14").Value2 = 2.3

# Row 15
$ws.Range("A15").Value2 = 'A 34302-2024'
$ws.Range("B15").Value2 = 45524
$ws.Range("C15").Value2 = 46081
$ws.Range("G15").Value2 = 1.9

# Row 16
$ws.Range("A16").Value2 = 'A 60803-2023'
$ws.Range("B16").Value2 = 45260
$ws.Range("C16").Value2 = 46081
$ws.Range("F16").Value2 = ""
$ws.Range("G16").Value2 = 1.6

# Row 17
$ws.Range("A17").Value2 = 'A 32596-2024'
$ws.Range("B17").Value2 = 45513.61667824074
$ws.Range("C17").Value2 = 46081
$ws.Range("G17").Value2 = 2.6

# Row 18
$ws.Range("A18").Value2 = 'A 49536-2025'
$ws.Range("B18").Value2 = 45939.4221875
$ws.Range("C18").Value2 = 46081
$ws.Range("G18").Value2 = 1.5

# Row 19
$ws.Range("A19").Value2 = 'A 49543-2025'
$ws.Range("B19").Value2 = 45939.42862268518
$ws.Range("C19").Value2 = 46081
$ws.Range("G19").Value2 = 1.4

# Row 20
$ws.Range("A20").Value2 = 'A 38631-2023'
$ws.Range("B20").Value2 = 45162
$ws.Range("C20").Value2 = 46081
$ws.Range("G20").Value2 = 0.8

# Row 21
$ws.Range("A21").Value2 = 'A 18090-2022'
$ws.Range("B21").Value2 = 44684
$ws.Range("C21").Value2 = 46081
$ws.Range("G21").Value2 = 4.9

# Row 22
$ws.Range("A22").Value2 = 'A 49549-2025'
$ws.Range("B22").Value2 = 45939
$ws.Range("C22").Value2 = 46081
$ws.Range("G22").Value2 = 0.5

# Row 23
$ws.Range("A23").Value2 = 'A 22195-2023'
$ws.Range("B23").Value2 = 45069.74605324074
$ws.Range("C23").Value2 = 46081
$ws.Range("F23").Value2 = 'Övriga Aktiebolag'
$ws.Range("G23").Value2 = 1.1
